# Actualización automática 2025-07-31 08:55:10
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M2").Value = 2820.46
$ws1.Range("I9").Value = 1697.4
$ws1.Range("M9").Value = 593.05
$ws1.Range("D12").Value = 915.84
$ws1.Range("M12").Value = 7519.12
$ws1.Range("M27").Value = 771.59
$ws1.Range("D58").Value = "3 de 56"
$ws1.Range("M58").Value = "11 de 56"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 2820.46
$ws2.Range("F9").Value = 2938.7
$ws2.Range("F12").Value = 8434.959999999999
$ws2.Range("F27").Value = 1212.58
$ws2.Range("F58").Value = 36289.64

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 4376.26
$ws3.Range("E3").Value = 615.9232000000002
$ws3.Range("F3").Value = 0.8766224765148842

$ws3.Range("D8").Value = 2010.3
$ws3.Range("E8").Value = -1260.3
$ws3.Range("F8").Value = 2.6804

$ws3.Range("D16").Value = 26362.77
$ws3.Range("E16").Value = 14027.4
$ws3.Range("F16").Value = 0.6527026254160356

$ws3.Range("D19").Value = 36351.22000000001
$ws3.Range("E19").Value = 19058.48560036207
$ws3.Range("F19").Value = 0.6560442724994821

$wb.Save()
